$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the crypto price refresh diff.
# Numeric-looking text values (e.g. "6.92", "0.999") need to be forced to
# text so Excel does not silently coerce them into real numbers (which would
# drop trailing zeros / reformat). We flip NumberFormat to "@" long enough to
# write the literal text, then ClearFormats() so no stray style id gets baked
# into the saved cell (these cells carry no explicit style in the original).

$ws.Range("D2").Value = '64.015.35'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '3.227.10'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '3.222.73'
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.33'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.75'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.27%  '
$ws.Range("D15").Value = '3.760.16'
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").Value = '3.229.59'
$ws.Range("E17").Value = '  -1.21%  '
$ws.Range("D18").Value = '64.019.64'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.62'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '468.81'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.12'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.709'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.76'
$ws.Range("D23").ClearFormats()
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.93'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("E35").Value = '  -2.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.95'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").Value = '0.0₃0745'
$ws.Range("E37").Value = '  +3.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '51.84'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0398'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '408.82'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.18'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  -2.33%  '
$ws.Range("D44").Value = '2.855.50'
$ws.Range("E44").Value = '  -7.16%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '128.88'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '36.16'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.93'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.112'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.19%  '
